$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue 'D2' '64.387.12'
Set-TextValue 'E2' '  -3.02%  '
Set-TextValue 'D3' '3.174.23'
Set-TextValue 'E3' '  -4.56%  '
Set-TextValue 'E4' '  -0.01%  '
Set-TextValue 'D5' '569.33'
Set-TextValue 'E5' '  -2.85%  '
Set-TextValue 'D6' '168.79'
Set-TextValue 'E6' '  -8.03%  '
Set-TextValue 'D7' '0.609'
Set-TextValue 'E7' '  -5.81%  '
Set-TextValue 'E8' '  -0.01%  '
Set-TextValue 'D9' '3.172.38'
Set-TextValue 'E10' '  -4.86%  '
Set-TextValue 'D11' '6.77'
Set-TextValue 'D12' '0.386'
Set-TextValue 'E12' '  -4.15%  '
Set-TextValue 'D13' '3.725.56'
Set-TextValue 'E13' '  -4.52%  '
Set-TextValue 'E14' '  -2.21%  '
Set-TextValue 'D15' '64.422.83'
Set-TextValue 'E15' '  -3.05%  '
Set-TextValue 'D16' '25.39'
Set-TextValue 'E16' '  -3.71%  '
Set-TextValue 'D17' '0.0000159'
Set-TextValue 'E17' '  -3.35%  '
Set-TextValue 'D18' '3.169.09'
Set-TextValue 'E18' '  -5.55%  '
Set-TextValue 'D19' '418.55'
Set-TextValue 'E19' '  -2.86%  '
Set-TextValue 'E20' '  -3.43%  '
Set-TextValue 'D21' '5.36'
Set-TextValue 'E21' '  -3.20%  '
Set-TextValue 'D22' '7.05'
Set-TextValue 'E22' '  -5.16%  '
Set-TextValue 'E23' '  -0.07%  '
Set-TextValue 'D24' '70.12'
Set-TextValue 'E24' '  -3.02%  '
Set-TextValue 'E25' '  +1.94%  '
Set-TextValue 'D26' '0.487'
Set-TextValue 'E26' '  -5.74%  '
Set-TextValue 'E27' '  -8.09%  '
Set-TextValue 'D28' '8.79'
Set-TextValue 'E28' '  -2.49%  '
Set-TextValue 'D29' '1.00'
Set-TextValue 'E29' '  +0.25%  '
Set-TextValue 'E30' '  -6.09%  '
Set-TextValue 'D31' '21.70'
Set-TextValue 'E31' '  -3.17%  '
Set-TextValue 'D33' '5.02'
Set-TextValue 'E33' '  -3.59%  '
Set-TextValue 'D34' '6.31'
Set-TextValue 'E34' '  -4.84%  '
Set-TextValue 'E35' '  -4.75%  '
Set-TextValue 'D36' '157.32'
Set-TextValue 'E36' '  -1.62%  '
Set-TextValue 'E37' '  -6.64%  '
Set-TextValue 'D38' '2.723.00'
Set-TextValue 'E39' '  -7.05%  '
Set-TextValue 'D40' '24.27'
Set-TextValue 'E40' '  -8.90%  '
Set-TextValue 'E41' '  -4.12%  '
Set-TextValue 'D42' '39.12'
Set-TextValue 'D43' '0.709'
Set-TextValue 'E43' '  -7.70%  '
Set-TextValue 'D44' '0.0620'
Set-TextValue 'E44' '  -7.00%  '
Set-TextValue 'D45' '5.61'
Set-TextValue 'E45' '  -6.55%  '
Set-TextValue 'E46' '  -3.97%  '
Set-TextValue 'D47' '21.60'
Set-TextValue 'E47' '  -7.98%  '
Set-TextValue 'D48' '293.24'
Set-TextValue 'E48' '  -7.59%  '
Set-TextValue 'E49' '  -0.03%  '
Set-TextValue 'D50' '2.00'
Set-TextValue 'E50' '  -13.82%  '
Set-TextValue 'D51' '0.0988'
